$wb = $excel.ActiveWorkbook

# The handoff xliff files were (re)generated for the
# c5913792-d698-4542-a08a-dd07ea4fac0d.md source file, so its "Latest
# Handoff Datetime" timestamp is refreshed on each per-language sheet, and
# the corresponding "Latest HO Xliff Generate Date" on the Overview sheet
# is refreshed to match the de-de value.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 6 on every sheet corresponds to c5913792-d698-4542-a08a-dd07ea4fac0d.md
$wsZhCn.Range("H6").Value = "2016-08-23 12:41:48"
$wsDeDe.Range("H6").Value = "2016-08-23 12:41:53"
$wsOverview.Range("G6").Value = "2016-08-23 12:41:53"
